$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 71, shifting existing rows 71..182 down to 72..183.
$ws.Rows(71).Insert()

# Populate the newly inserted row 71 with the new data record.
$ws.Cells.Item(71, 1).Value = 4
$ws.Cells.Item(71, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(71, 3).Value = "Los Lagos"
$ws.Cells.Item(71, 4).Value = 45175
$ws.Cells.Item(71, 5).Value = 10
$ws.Cells.Item(71, 6).Value = "Fruta"
$ws.Cells.Item(71, 7).Value = 100104
$ws.Cells.Item(71, 8).Value = "Frutos de pepita"
$ws.Cells.Item(71, 9).Value = 100104003
$ws.Cells.Item(71, 10).Value = "Membrillo"
$ws.Cells.Item(71, 11).Value = "Champion"
$ws.Cells.Item(71, 12).Value = "Primera"
$ws.Cells.Item(71, 13).Value = 40
$ws.Cells.Item(71, 14).Value = 16000
$ws.Cells.Item(71, 15).Value = 16000
$ws.Cells.Item(71, 16).Value = 16000
$ws.Cells.Item(71, 17).Value = "$/caja 18 kilos empedrada"
$ws.Cells.Item(71, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(71, 19).Value = 889
$ws.Cells.Item(71, 20).Value = 18
